# Daily update at 8 AM UTC
# Append the next day's row (46) to the "Wins Over Time" tracking sheet and
# restore the previous last row (45) to the regular datetime format, since it
# is no longer the final row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45 was previously the last row (date-only format). It is no longer the
# last row, so give it back the standard datetime number format used by the
# rest of the data rows.
$ws.Range("A45").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 46 becomes the last row: give it the date-only number format and
# fill in the day's values.
$ws.Range("A46").Value = 45786
$ws.Range("A46").NumberFormat = "YYYY-MM-DD"
$ws.Range("B46").Value = 186
$ws.Range("C46").Value = 198
$ws.Range("D46").Value = 193
